$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (before current row 124) to make
# room for the newest week's measurements. Existing rows 124:137 shift down to 126:139,
# which matches every other observed change in the diff (each row's content is simply
# the prior row's content, now one "slot" further down).
$ws.Rows("124:125").Insert()

# Row 124 — Apio, Primera, Región de Coquimbo, new week (2021-09-10 / serial 44449)
$ws.Cells.Item(124, 1).Value = 11
$ws.Cells.Item(124, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(124, 3).Value = "Bíobío"
$ws.Cells.Item(124, 4).Value2 = 44449
$ws.Cells.Item(124, 5).Value = 8
$ws.Cells.Item(124, 6).Value = 100112017
$ws.Cells.Item(124, 7).Value = "Apio"
$ws.Cells.Item(124, 8).Value = "Americana (o)"
$ws.Cells.Item(124, 9).Value = "Primera"
$ws.Cells.Item(124, 10).Value = 100
$ws.Cells.Item(124, 11).Value = 8500
$ws.Cells.Item(124, 12).Value = 9000
$ws.Cells.Item(124, 13).Value = 8750
$ws.Cells.Item(124, 14).Value = "$/docena de matas"
$ws.Cells.Item(124, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(124, 16).Value = 1458
$ws.Cells.Item(124, 17).Value = 6
$ws.Cells.Item(124, 18).Value = "Hortaliza"

# Row 125 — Apio, Segunda, Región de Coquimbo, same new week
$ws.Cells.Item(125, 1).Value = 11
$ws.Cells.Item(125, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(125, 3).Value = "Bíobío"
$ws.Cells.Item(125, 4).Value2 = 44449
$ws.Cells.Item(125, 5).Value = 8
$ws.Cells.Item(125, 6).Value = 100112017
$ws.Cells.Item(125, 7).Value = "Apio"
$ws.Cells.Item(125, 8).Value = "Americana (o)"
$ws.Cells.Item(125, 9).Value = "Segunda"
$ws.Cells.Item(125, 10).Value = 50
$ws.Cells.Item(125, 11).Value = 8000
$ws.Cells.Item(125, 12).Value = 8000
$ws.Cells.Item(125, 13).Value = 8000
$ws.Cells.Item(125, 14).Value = "$/docena de matas"
$ws.Cells.Item(125, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(125, 16).Value = 1333
$ws.Cells.Item(125, 17).Value = 6
$ws.Cells.Item(125, 18).Value = "Hortaliza"
